$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the names from column C into column B (same values/shared strings)
$ws.Range("B4:B32").Value = $ws.Range("C4:C32").Value()

# Seed cell for the running counter
$ws.Range("A3").Value = 0

# Running count formulas down column A
for ($r = 4; $r -le 33; $r++) {
    $prev = $r - 1
    $ws.Range("A$r").Formula = "=IF(B$r=C$r,1+A$prev,""____________"")"
}

$ws.Range("B27").Select()
